$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price table (Price / Volume(1h) columns, plus a few
# re-ranked coin rows) with the latest scrape. All cells in B:E are stored
# as literal text in the source sheet, so price strings that happen to look
# like plain decimals (e.g. "310.03", "1.80") are written with a leading
# apostrophe -- exactly like typing ' in Excel -- to force literal text
# entry instead of letting them be auto-converted to Number cells (which
# would also silently drop significant trailing zeros, e.g. "1.80" -> 1.8).

$ws.Range("D2").Value = "40.607.76"
$ws.Range("E2").Value = "  -2.86%  "

$ws.Range("D3").Value = "2.368.81"
$ws.Range("E3").Value = "  -4.40%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'310.03"
$ws.Range("E5").Value = "  -2.83%  "

$ws.Range("D6").Value = "'86.89"
$ws.Range("E6").Value = "  -6.92%  "

$ws.Range("D7").Value = "'0.527"
$ws.Range("E7").Value = "  -4.86%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "  -5.21%  "

$ws.Range("D10").Value = "'0.0842"
$ws.Range("E10").Value = "  -4.58%  "

$ws.Range("D11").Value = "'30.48"
$ws.Range("E11").Value = "  -8.36%  "

$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").Value = "2.734.21"
$ws.Range("E13").Value = "  -4.35%  "

$ws.Range("D14").Value = "'6.53"
$ws.Range("E14").Value = "  -5.89%  "

$ws.Range("D15").Value = "'14.93"
$ws.Range("E15").Value = "  -4.75%  "

$ws.Range("D16").Value = "2.378.67"
$ws.Range("E16").Value = "  -4.00%  "

$ws.Range("D17").Value = "'0.755"
$ws.Range("E17").Value = "  -6.27%  "

$ws.Range("D18").Value = "40.545.68"
$ws.Range("E18").Value = "  -2.91%  "

$ws.Range("D19").Value = "0.0₃0905"
$ws.Range("E19").Value = "  -5.03%  "

$ws.Range("D20").Value = "'6.11"
$ws.Range("E20").Value = "  -6.02%  "

$ws.Range("D21").Value = "'68.59"
$ws.Range("E21").Value = "  -3.76%  "

$ws.Range("D22").Value = "'10.72"
$ws.Range("E22").Value = "  -5.90%  "

$ws.Range("D23").Value = "'232.92"
$ws.Range("E23").Value = "  -3.75%  "

$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "  -4.75%  "

$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").Value = "'1.80"
$ws.Range("E26").Value = "  -8.10%  "

$ws.Range("D27").Value = "'23.58"
$ws.Range("E27").Value = "  -6.57%  "

$ws.Range("D28").Value = "'2.19"
$ws.Range("E28").Value = "  -3.18%  "

$ws.Range("D29").Value = "'9.26"
$ws.Range("E29").Value = "  -5.28%  "

$ws.Range("D30").Value = "'33.69"
$ws.Range("E30").Value = "  -9.00%  "

$ws.Range("D31").Value = "'152.37"
$ws.Range("E31").Value = "  -4.12%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").Value = "'5.19"
$ws.Range("E33").Value = "  -6.23%  "

$ws.Range("D34").Value = "'0.0726"
$ws.Range("E34").Value = "  -5.21%  "

$ws.Range("E35").Value = "  -4.91%  "

$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  -2.49%  "

$ws.Range("D37").Value = "'2.74"
$ws.Range("E37").Value = "  -6.64%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.0993"
$ws.Range("E38").Value = "  -5.35%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.71"
$ws.Range("E39").Value = "  -9.80%  "

$ws.Range("D40").Value = "'1.69"
$ws.Range("E40").Value = "  -9.82%  "

$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").Value = "'2.41"
$ws.Range("E41").Value = "  -4.57%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'3.82"
$ws.Range("E42").Value = "  -5.28%  "

$ws.Range("D43").Value = "1.951.27"
$ws.Range("E43").Value = "  -2.64%  "

$ws.Range("D44").Value = "'0.0269"
$ws.Range("E44").Value = "  -5.84%  "

$ws.Range("D45").Value = "'17.80"
$ws.Range("E45").Value = "  -7.47%  "

$ws.Range("D46").Value = "'9.42"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").Value = "'2.69"
$ws.Range("E47").Value = "  -10.04%  "

$ws.Range("D48").Value = "2.599.19"
$ws.Range("E48").Value = "  -4.30%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'72.26"
$ws.Range("E49").Value = "  -5.91%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'92.68"
$ws.Range("E50").Value = "  -5.77%  "

$ws.Range("D51").Value = "'50.30"
$ws.Range("E51").Value = "  -4.65%  "
